$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a row below "Decsion Tree" (row 7) for "Neureal Network" first, so
# that it claims the next shared-string slot before "kNN" does.
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Neureal Network"
$ws.Range("B8").Value = 82.4

# Insert a row above "Decsion Tree" (row 7) for "kNN", pushing
# "Decsion Tree", "Neureal Network" and "Random Forest Classifier" down.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "kNN"
$ws.Range("B7").Value = 80.1

# Re-apply the sort over the now-larger data range so the sheet's sort
# state tracks the expanded A4:B10 block.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B4:B10"))
$sortObj.SetRange($ws.Range("A4:B10"))
$sortObj.Header = 0
$sortObj.Apply()

# Grow the bar chart's category/value series so it plots the two new rows
# as well (was Sheet1!$A$4:$A$8 / $B$4:$B$8, now extends to row 10).
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Sheet1!`$B`$3,Sheet1!`$A`$4:`$A`$10,Sheet1!`$B`$4:`$B`$10,1)"

# Re-select per the target sheet view and bump the zoom level.
$ws.Range("H1").Select()
$excel.ActiveWindow.Zoom = 150
